$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in rows 4-7 ---
$ws.Range("B4").Value = -0.00000003912214467726298
$ws.Range("C4").Value = 2664.798010691252
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = 35238095.23809128
$ws.Range("C5").Value = 4630783.093411878
$ws.Range("D5").Value = 0

$ws.Range("C6").Value = -0.0000002522017896953197

$ws.Range("C7").Value = 0.000000002404676588544419

# --- Add new rows 13-15 ---
$ws.Range("A13").Value = "Electrification + Bio-based feedstock"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = -0.0000008923596138608058
$ws.Range("D13").Value = 0

$ws.Range("A14").Value = "Conventional + Bio-based feedstock"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 9428.791630345906
$ws.Range("D14").Value = 0.0000002114406099842145

$ws.Range("A15").Value = "Conventional + Bio-based feedstock with CC"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 16384990.04349792
$ws.Range("D15").Value = 28007172.12752622

# Apply the same style as the other label cells in column A (A4:A12) to the new A13:A15 cells
$ws.Range("A4").Copy()
$ws.Range("A13:A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
